$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "G2"
$ws.Range("B3").Value = "Mask1"
$ws.Range("C3").Value = "Daily"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 45860
$ws.Range("E3").NumberFormat = "YYYY-MM-DD"
$ws.Range("F3").Value = 30
